$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.19678008556366
$ws.Range("B1").Value = 2.482578277587891
$ws.Range("C1").Value = 4.18312931060791
$ws.Range("D1").Value = 2.076373100280762
$ws.Range("E1").Value = 1.183169603347778
